$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" (column D) values look like plain numbers (e.g. "1.00", "522.98").
# Excel would normally auto-convert such text to a numeric value when assigned via
# .Value, but the source data stores them as text. Temporarily force the whole
# D2:D51 column to Text format, assign the new values, then clear the formatting
# again so the cells end up with no explicit style (matching the original layout)
# while keeping their content as text.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '58.930.03'
$ws.Range("E2").Value = '  +3.14%  '
$ws.Range("D3").Value = '3.094.26'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '522.98'
$ws.Range("E5").Value = '  +2.08%  '
$ws.Range("D6").Value = '143.67'
$ws.Range("E6").Value = '  +1.00%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.440'
$ws.Range("E8").Value = '  +1.00%  '
$ws.Range("E9").Value = '  +1.01%  '
$ws.Range("E10").Value = '  +1.41%  '
$ws.Range("E11").Value = '  +2.90%  '
$ws.Range("D12").Value = '3.626.08'
$ws.Range("E12").Value = '  +0.98%  '
$ws.Range("E13").Value = '  +1.02%  '
$ws.Range("D14").Value = '26.82'
$ws.Range("E14").Value = '  +4.67%  '
$ws.Range("D15").Value = '0.0000167'
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("D16").Value = '58.913.24'
$ws.Range("E16").Value = '  +2.80%  '
$ws.Range("D17").Value = '3.099.22'
$ws.Range("E17").Value = '  +1.38%  '
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("D19").Value = '12.94'
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '8.11'
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").Value = '343.49'
$ws.Range("E21").Value = '  +2.49%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '0.507'
$ws.Range("E23").Value = '  +1.84%  '
$ws.Range("D24").Value = '65.70'
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("E28").Value = '  +4.77%  '
$ws.Range("D29").Value = '7.25'
$ws.Range("E29").Value = '  +2.95%  '
$ws.Range("E30").Value = '  +2.44%  '
$ws.Range("D31").Value = '1.20'
$ws.Range("E31").Value = '  +3.08%  '
$ws.Range("D32").Value = '20.99'
$ws.Range("E32").Value = '  +1.43%  '
$ws.Range("D33").Value = '154.83'
$ws.Range("E33").Value = '  +0.42%  '
$ws.Range("E34").Value = '  +2.98%  '
$ws.Range("D35").Value = '6.14'
$ws.Range("E35").Value = '  +4.49%  '
$ws.Range("D36").Value = '26.83'
$ws.Range("E36").Value = '  +2.59%  '
$ws.Range("E37").Value = '  +4.75%  '
$ws.Range("D38").Value = '0.0685'
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").Value = '3.93'
$ws.Range("E39").Value = '  +2.53%  '
$ws.Range("D40").Value = '3.135.04'
$ws.Range("E40").Value = '  +1.06%  '
$ws.Range("D41").Value = '36.76'
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").Value = '0.665'
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("E44").Value = '  +5.80%  '
$ws.Range("D45").Value = '2.286.34'
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("D46").Value = '0.0255'
$ws.Range("E46").Value = '  +1.30%  '
$ws.Range("D47").Value = '20.92'
$ws.Range("E47").Value = '  +3.35%  '
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("E49").Value = '  +2.66%  '
$ws.Range("D50").Value = '0.754'
$ws.Range("E50").Value = '  +9.41%  '
$ws.Range("D51").Value = '261.77'
$ws.Range("E51").Value = '  +10.82%  '

$priceRange.ClearFormats()
